$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ row = 17; ts = "2025-08-20 23:15:23" },
    @{ row = 18; ts = "2025-08-20 23:15:38" }
)

foreach ($r in $rows) {
    $i = $r.row
    $ws.Cells.Item($i, 1).Value = "test"
    $ws.Cells.Item($i, 2).Value = 1934
    $ws.Cells.Item($i, 3).Value = 1661
    $ws.Cells.Item($i, 4).Value = 4
    $ws.Cells.Item($i, 5).Value = 1633
    $ws.Cells.Item($i, 6).NumberFormat = "@"
    $ws.Cells.Item($i, 6).Value = "36.96"
    $ws.Cells.Item($i, 6).Style = "Normal"
    $ws.Cells.Item($i, 7).NumberFormat = "@"
    $ws.Cells.Item($i, 7).Value = "31.75"
    $ws.Cells.Item($i, 7).Style = "Normal"
    $ws.Cells.Item($i, 8).NumberFormat = "@"
    $ws.Cells.Item($i, 8).Value = "0.08"
    $ws.Cells.Item($i, 8).Style = "Normal"
    $ws.Cells.Item($i, 9).NumberFormat = "@"
    $ws.Cells.Item($i, 9).Value = "31.21"
    $ws.Cells.Item($i, 9).Style = "Normal"
    $ws.Cells.Item($i, 10).Value = $r.ts
}
